# smart city studio 1.1
# Add a trailing colon after each class-meeting date in the schedule
# table, and turn the "Apr 13" row's Workshops/Critique cells into the
# Midpoint Critique milestone.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Column 2 of rows 2..18 holds the class-meeting dates; append ":" to
# each one in place (row 1 is the "Dates" header, left untouched).
# Cell.Range.Text includes the trailing cell-end marks (CR + cell
# marker), so strip those two characters before appending the colon.
for ($r = 2; $r -le 18; $r++) {
    $cell = $t.Cell($r, 2)
    $txt = $cell.Range.Text
    $dateOnly = $txt.Substring(0, $txt.Length - 2)
    $cell.Range.Text = $dateOnly + ":"
}

# Update the Apr 13 (row 13) milestone cells: Workshops -> Midpoint
# Critique, Critique -> Project prototype.
$t.Cell(13, 3).Range.Text = "Midpoint Critique"
$t.Cell(13, 4).Range.Text = "Project prototype"
